$d = $word.ActiveDocument

$p3 = $d.Paragraphs(3)
$p3.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="00E8765B" w:rsidR="00E8765B" w:rsidRDefault="00E8765B"><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t/></w:r><w:r w:rsidR="2474C2E1419B439FA1711D54B1990E30"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="2474C2E1419B439FA1711D54B1990E30"><w:instrText xml:space="preserve"> REF Art1 \h </w:instrText></w:r><w:r w:rsidR="2474C2E1419B439FA1711D54B1990E30"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="2474C2E1419B439FA1711D54B1990E30"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>Artifact1</w:t></w:r><w:r w:rsidR="2474C2E1419B439FA1711D54B1990E30"><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$p5 = $d.Paragraphs(5)
$p5.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979"><w:r w:rsidR="00E61FB8"><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t/></w:r><w:bookmarkStart w:name="Art1" w:id="96905454422089742970739521747084098851"/><w:r w:rsidR="00E61FB8"><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>Definition of Artifact1</w:t></w:r><w:bookmarkEnd w:id="96905454422089742970739521747084098851"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$p7 = $d.Paragraphs(7)
$p7.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="00F5495F" w:rsidR="00730F00" w:rsidRDefault="00730F00"><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t/></w:r><w:r w:rsidR="2474C2E1419B439FA1711D54B1990E30"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="2474C2E1419B439FA1711D54B1990E30"><w:instrText xml:space="preserve"> REF Art1 \h </w:instrText></w:r><w:r w:rsidR="2474C2E1419B439FA1711D54B1990E30"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="2474C2E1419B439FA1711D54B1990E30"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>Artifact1</w:t></w:r><w:r w:rsidR="2474C2E1419B439FA1711D54B1990E30"><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
